# Update cryptos list worksheet with the latest scraped price/volume data.
# Several numeric-looking "Price" values (column D) must remain plain text
# (they use "." as a thousands marker, e.g. "600.54"), so for those cells we
# temporarily force a text number format before writing the value and then
# restore the cell's original style to avoid leaving stray formatting behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.735.59"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.677.35"
$ws.Range("E3").Value = "  -0.83%  "
$ws.Range("E4").Value = "  -0.01%  "
$c = $ws.Range("D5")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "600.54"
$c.Style = $origStyle
$ws.Range("E5").Value = "  -0.89%  "
$c = $ws.Range("D6")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "156.87"
$c.Style = $origStyle
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +6.06%  "
$ws.Range("E9").Value = "  +5.10%  "
$ws.Range("E10").Value = "  +0.01%  "
$ws.Range("E11").Value = "  -2.71%  "
$ws.Range("E12").Value = "  -0.12%  "
$ws.Range("E13").Value = "  -2.49%  "
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").Value = "3.157.43"
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("D16").Value = "65.609.23"
$ws.Range("E16").Value = "  +0.00%  "
$ws.Range("D17").Value = "2.675.43"
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("E20").Value = "  +1.38%  "
$c = $ws.Range("D21")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "352.35"
$c.Style = $origStyle
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("E22").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +5.61%  "
$ws.Range("E25").Value = "  -1.47%  "
$c = $ws.Range("D26")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "1.64"
$c.Style = $origStyle
$ws.Range("E26").Value = "  +0.20%  "
$ws.Range("E27").Value = "  -2.22%  "
$ws.Range("E28").Value = "  -5.77%  "
$c = $ws.Range("D29")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "8.07"
$c.Style = $origStyle
$ws.Range("E29").Value = "  -3.17%  "
$ws.Range("E31").Value = "  -2.45%  "
$c = $ws.Range("D32")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "528.41"
$c.Style = $origStyle
$ws.Range("E32").Value = "  -1.97%  "
$ws.Range("E33").Value = "  -1.44%  "
$c = $ws.Range("D34")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "6.49"
$c.Style = $origStyle
$ws.Range("E34").Value = "  -3.60%  "
$ws.Range("B36").Value = "PolygonEcosystemToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$c = $ws.Range("D36")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "0.424"
$c.Style = $origStyle
$ws.Range("E36").Value = "  -1.69%  "
$ws.Range("B37").Value = "EthereumClassic"
$ws.Range("C37").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D37")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "20.67"
$c.Style = $origStyle
$ws.Range("E37").Value = "  -0.23%  "
$ws.Range("E38").Value = "  -0.02%  "
$c = $ws.Range("D39")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "157.97"
$c.Style = $origStyle
$ws.Range("E39").Value = "  -2.80%  "
$ws.Range("E40").Value = "  -2.37%  "
$c = $ws.Range("D42")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "164.21"
$c.Style = $origStyle
$ws.Range("E42").Value = "  -2.16%  "
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("E44").Value = "  +2.89%  "
$ws.Range("E45").Value = "  -0.23%  "
$c = $ws.Range("D46")
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = "22.92"
$c.Style = $origStyle
$ws.Range("E46").Value = "  -2.73%  "
$ws.Range("E47").Value = "  +17.52%  "
$ws.Range("E48").Value = "  -2.67%  "
$ws.Range("E50").Value = "  +2.33%  "
$ws.Range("E51").Value = "  -4.88%  "
